# Append new substance rows to the data table (Tabelle1), matching the
# four new entries added by the author: a ligand (HMTETA) and three
# macroinitiator batches (HW45/49/55-PtBMA...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: HMTETA (ligand, liquid)
$ws.Range("A12").Value = "HMTETA"
$ws.Range("B12").Value = "ligand"
$ws.Range("C12").Value = 230.39
$ws.Range("D12").Value = 0.847
$ws.Range("E12").Value = "liquid"

# Row 13: HW45-PtBMA (makroinitiator, solid)
$ws.Range("A13").Value = "HW45-PtBMA"
$ws.Range("B13").Value = "makroinitiator"
$ws.Range("C13").Value = 8156
$ws.Range("E13").Value = "solid"

# Row 14: HW49-PtBMA-END (makroinitiator, solid)
$ws.Range("A14").Value = "HW49-PtBMA-END"
$ws.Range("B14").Value = "makroinitiator"
$ws.Range("C14").Value = 8500
$ws.Range("E14").Value = "solid"

# Row 15: HW55-PtBMA-Prod (makroinitiator, solid)
$ws.Range("A15").Value = "HW55-PtBMA-Prod"
$ws.Range("B15").Value = "makroinitiator"
$ws.Range("C15").Value = 8674
$ws.Range("E15").Value = "solid"

# Widen column A so the longer substance names fit (as seen by the
# author re-sizing it manually before saving).
$ws.Columns.Item(1).ColumnWidth = 17.333333333333332

# Leave selection on the last touched cell, like the saved file shows.
$ws.Range("E15").Select() | Out-Null
